$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 160, shifting rows 160:165 down to 161:166
$ws.Rows("160").Insert()

$ws.Cells.Item(160, 1).Value = 10
$ws.Cells.Item(160, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(160, 3).Value = "La Araucanía"
$ws.Cells.Item(160, 4).Value = 44585
$ws.Cells.Item(160, 5).Value = 9
$ws.Cells.Item(160, 6).Value = "Fruta"
$ws.Cells.Item(160, 7).Value = 100103
$ws.Cells.Item(160, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(160, 9).Value = 100103002
$ws.Cells.Item(160, 10).Value = "Ciruela"
$ws.Cells.Item(160, 11).Value = "Black Amber"
$ws.Cells.Item(160, 12).Value = "Primera"
$ws.Cells.Item(160, 13).Value = 300
$ws.Cells.Item(160, 14).Value = 12000
$ws.Cells.Item(160, 15).Value = 12000
$ws.Cells.Item(160, 16).Value = 12000
$ws.Cells.Item(160, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(160, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(160, 19).Value = 667
$ws.Cells.Item(160, 20).Value = 18
